$d = $word.ActiveDocument

# 1. Header: merge the split "Ge" / "offrey A. Reed" runs into a single
#    run containing "Geoffrey A. Reed" (same text, same formatting -
#    just one run instead of two). MatchCase must stay $false here
#    because the header text is displayed as caps (w:caps) even though
#    the underlying characters are mixed-case.
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)
    $hdr.Range.Find.Execute("Geoffrey A. Reed", $false, $false, $false, $false, $false, $true, 1, $false, "Geoffrey A. Reed", 2) | Out-Null
}

# 2. Style "Footnote Reference": bump the run font size from 8pt (16
#    half-points) to 10pt (20 half-points).
$footnoteRefStyle = $d.Styles.Item("Footnote Reference")
$footnoteRefStyle.Font.Size = 10

Write-Output "done"
